$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Selector" column (C) is dropped; the old "Title" column (D) takes its
# place, the xpath selector value that lived in C6 moves into the new D
# column, and two fresh columns record the mobile + web testcase results
# (Chromium / Pixel9Pro_API35).

# Column C becomes "Title" (was column D).
$ws.Range("C1").Value = "Title"
$ws.Range("C2").Value = "Login | TT-Planer"
$ws.Range("C3").Value = "Übersicht | TT-Planer"
$ws.Range("C4").Value = "Passwort vergessen | TT-Planer"
$ws.Range("C5").Value = "Mein Profil | TT-Planer"
$ws.Range("C6").Value = "Abwesenheiten | TT-Planer"

# Column D is repurposed: header becomes "Chromium" and only keeps the
# xpath selector value (moved over from the old column C) on row 6.
$ws.Range("D1").Value = "Chromium"
$ws.Range("D2").ClearContents()
$ws.Range("D3").ClearContents()
$ws.Range("D4").ClearContents()
$ws.Range("D5").ClearContents()
$ws.Range("D6").Value = "//div[@id='createAbsenceModal']"

# New column E: header for the mobile test device.
$ws.Range("E1").Value = "Pixel9Pro_API35"
$ws.Range("E1").NumberFormat = "@"

# Restore bestFit-style widths for the now-shifted columns and size the new one
# (values chosen so the engine's column-width quantisation lands as close as
# possible to the original bestFit widths of 26.77734375 / 28.88671875 / 20.6640625).
$ws.Columns.Item(3).ColumnWidth = 26
$ws.Columns.Item(4).ColumnWidth = 28
$ws.Columns.Item(5).ColumnWidth = 19.8333333333333

$ws.Range("F8").Select()
